$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "1.005") are stored as literal text, matching the source data,
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.319.70"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "1.706.12"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "223.71"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.2657"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").Value = "0.06575"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "20.76"
$ws.Range("E10").Value = "  -4.61%  "
$ws.Range("D11").Value = "0.07624"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "4.517"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "1.722.94"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "1.941.44"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "0.5760"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "0.0₅8132"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "67.53"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "27.321.99"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "215.26"
$ws.Range("E19").Value = "  -4.41%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").Value = "10.38"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "5.939"
$ws.Range("E23").Value = "  -3.22%  "
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "143.55"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("D26").Value = "1.716"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").Value = "0.1203"
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").Value = "7.209"
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("D29").Value = "16.05"
$ws.Range("E29").Value = "  -4.73%  "
$ws.Range("D30").Value = "0.05366"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "3.472"
$ws.Range("E32").Value = "  -3.53%  "
$ws.Range("D33").Value = "3.401"
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("D34").Value = "1.639"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").Value = "2.878"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").Value = "2.416"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "0.9447"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").Value = "0.5797"
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").Value = "0.01625"
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("D40").Value = "5.760"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "1.038.31"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").Value = "101.09"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "1.849.16"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").Value = "0.0₈117"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").Value = "57.66"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").Value = "0.4520"
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("D49").Value = "0.9979"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "8.058"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "0.05233"
$ws.Range("E51").Value = "  -1.05%  "

# Restore the default cell style on column D so only the values changed
# (the temporary Text number format above would otherwise leave a new
# style index applied to these cells).
$ws.Range("D2:D51").Style = "Normal"
